# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.953.41'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '1.911.53'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'" + '319.81'
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').Value = "'" + '0.9995'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = "'" + '0.5037'
$ws.Range('E7').Value = '  -2.45%  '
$ws.Range('D8').Value = "'" + '0.4038'
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('D10').Value = "'" + '1.102'
$ws.Range('E10').Value = '  -1.39%  '
$ws.Range('D11').Value = "'" + '41.94'
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('D12').Value = "'" + '24.14'
$ws.Range('E12').Value = '  +3.62%  '
$ws.Range('D13').Value = '1.914.91'
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').Value = "'" + '6.390'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = "'" + '7.210'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').Value = "'" + '1.001'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = "'" + '92.14'
$ws.Range('E17').Value = '  -2.78%  '
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('D19').Value = "'" + '0.06499'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').Value = "'" + '18.15'
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').Value = "'" + '1.000'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = "'" + '5.931'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = '29.986.67'
$ws.Range('E23').Value = '  -0.84%  '
$ws.Range('D24').Value = "'" + '11.30'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = "'" + '2.192'
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('D26').Value = "'" + '22.18'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').Value = '2.134.47'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = "'" + '162.24'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').Value = "'" + '2.299'
$ws.Range('E29').Value = '  -2.84%  '
$ws.Range('D30').Value = "'" + '128.92'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = "'" + '1.127'
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').Value = "'" + '0.1035'
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('D33').Value = "'" + '5.929'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('D34').Value = "'" + '3.819'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('D35').Value = "'" + '5.421'
$ws.Range('E35').Value = '  +3.12%  '
$ws.Range('D36').Value = "'" + '0.02440'
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('D37').Value = "'" + '0.06397'
$ws.Range('D38').Value = "'" + '0.2143'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = "'" + '1.194'
$ws.Range('E39').Value = '  -2.06%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'" + '8.702'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = "'" + '0.6464'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('D42').Value = "'" + '11.37'
$ws.Range('E42').Value = '  -3.69%  '
$ws.Range('D43').Value = "'" + '1.211'
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('D44').Value = "'" + '2.213'
$ws.Range('E44').Value = '  +7.57%  '
$ws.Range('D45').Value = "'" + '13.30'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').Value = "'" + '0.6033'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').Value = "'" + '3.633'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('D48').Value = "'" + '122.10'
$ws.Range('E48').Value = '  -1.84%  '
$ws.Range('D49').Value = "'" + '1.204'
$ws.Range('E49').Value = '  -2.57%  '
$ws.Range('D50').Value = "'" + '78.88'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = "'" + '1.128'
$ws.Range('E51').Value = '  -3.05%  '
